# Update the two sample-size caption textboxes that live inside the
# "Group 228" shape on slide 1 (bottom-right legend of the figure).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the group shape that contains both captions by name, regardless
# of its absolute position in the shape collection.
$grp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Group 228") {
        $grp = $candidate
    }
}

# --- "N>100 samples" textbox -> "Up to 1000 samples", widened box ---
$tb1 = $grp.GroupItems.Item("TextBox 246")
$tr1 = $tb1.TextFrame.TextRange

# Drop the ">"-run and the "100 samples"-run, leaving just the first
# run (which already carries dirty="0"), then rewrite its text.
$tr1.Characters(2, $tr1.Length - 1).Delete()
$tr1.Characters(1, 1).Text = "Up to 1000 samples"

# Widen the textbox to fit the new, longer caption.
$tb1.Width = 112.371815

# --- "N>1000 samples" textbox: merge the trailing two runs ---
$tb2 = $grp.GroupItems.Item("TextBox 262")
$tr2 = $tb2.TextFrame.TextRange

# Runs are: "N" / ">" / "1000 " / "samples". Remove the "1000 " run and
# fold its text into the "samples" run so the caption reads as a single
# "1000 samples" run (matching the other caption's single trailing run).
$tr2.Characters(3, 5).Delete()
$tr2.Characters(3, 7).Text = "1000 samples"
